# "tried new udi's math game background"
# Append 4 more copies of the last existing data row (row 102) to the
# Users sheet, producing new rows 103-106 with identical content
# (Username=moses, Password=bro, ID=1234, Email=m@g.c, Gender=Male, balance=0).
#
# We copy the existing row instead of writing cell values directly so that
# the new cells pick up the exact same cell types/styles as the source row
# (in particular the "ID" column is stored as text "1234" via the shared
# string table, not as a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRow = $ws.Range("A102:F102")

$sourceRow.Copy($ws.Range("A103:F103"))
$sourceRow.Copy($ws.Range("A104:F104"))
$sourceRow.Copy($ws.Range("A105:F105"))
$sourceRow.Copy($ws.Range("A106:F106"))
